# The header row (row 1) contains the "T2" condition labels. Strip the
# trailing "T2" suffix so the sheet matches the non-T2 condition naming
# used elsewhere in the workbook:
#   squareT2 -> square, loc1T2 -> loc1, loc2T2 -> loc2, corrAnsT2 -> corrAns

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "square"
$ws.Range("B1").Value = "loc1"
$ws.Range("C1").Value = "loc2"
$ws.Range("D1").Value = "corrAns"
